# Apply updated cryptocurrency price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.644.75'
$ws.Range('E2').Value = '  +3.59%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.699.23'
$ws.Range('E3').Value = '  +2.41%  '

$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.86'
$ws.Range('E5').Value = '  +2.98%  '

$ws.Range('E6').Value = '  -0.09%  '

$ws.Range('E7').Value = '  +1.64%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4041'
$ws.Range('E8').Value = '  +2.52%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.546'
$ws.Range('E9').Value = '  +9.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '54.94'
$ws.Range('E10').Value = '  +12.45%  '

$ws.Range('E11').Value = '  -0.07%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08812'
$ws.Range('E12').Value = '  +2.35%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.302'
$ws.Range('E13').Value = '  +9.24%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.48'
$ws.Range('E14').Value = '  +3.67%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001334'
$ws.Range('E15').Value = '  +2.42%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.649'
$ws.Range('E16').Value = '  +6.55%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.697.24'
$ws.Range('E17').Value = '  +1.89%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '101.40'
$ws.Range('E18').Value = '  +1.57%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07100'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.81'
$ws.Range('E20').Value = '  +4.73%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.907'
$ws.Range('E21').Value = '  +4.17%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.0000'
$ws.Range('E22').Value = '  -0.16%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.20'
$ws.Range('E23').Value = '  +3.35%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.636.46'
$ws.Range('E24').Value = '  +3.63%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.045'
$ws.Range('E25').Value = '  +12.10%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.334'
$ws.Range('E26').Value = '  +0.96%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.45'
$ws.Range('E27').Value = '  +3.62%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '159.90'

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.247'
$ws.Range('E29').Value = '  +1.25%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.19'
$ws.Range('E30').Value = '  +3.64%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.660'
$ws.Range('E31').Value = '  +19.12%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.883.79'
$ws.Range('E32').Value = '  +1.98%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.105'
$ws.Range('E33').Value = '  -1.94%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.413'
$ws.Range('E34').Value = '  +13.26%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08584'
$ws.Range('E35').Value = '  +0.31%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '11.36'
$ws.Range('E36').Value = '  +10.58%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2763'
$ws.Range('E37').Value = '  +5.21%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.951'
$ws.Range('E38').Value = '  -0.82%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.80'
$ws.Range('E39').Value = '  +3.41%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02807'
$ws.Range('E40').Value = '  +11.71%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.09084'
$ws.Range('E41').Value = '  +3.85%  '

$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7789'
$ws.Range('E42').Value = '  +3.57%  '

$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.467'
$ws.Range('E43').Value = '  +1.41%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7299'
$ws.Range('E44').Value = '  +4.37%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.62'
$ws.Range('E45').Value = '  +5.89%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.520'
$ws.Range('E46').Value = '  +6.52%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.224'
$ws.Range('E47').Value = '  +4.66%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.397'
$ws.Range('E48').Value = '  +21.61%  '

$ws.Range('E49').Value = '  -0.08%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '141.98'
$ws.Range('E50').Value = '  +2.02%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.08040'
$ws.Range('E51').Value = '  +3.88%  '
